$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-covering dates) ---
# A8 holds "Volume 29   Number  49" -> "...50"
$ws.Range("A8").Value = "Volume 29   Number  50"
# C9 holds "Report Covering the Week  12/5/2022  Through  12/11/2022"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# --- Row 15 (Murder): only N15 changes ---
$ws.Range("N15").Value = -72.727272727272

# --- Row 16 (Rape): D16/E16/G16/H16 go from "N/A"/"***.*" text to real numbers ---
$ws.Range("F17").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 2
$ws.Range("H17").Copy($ws.Range("E16"))
$ws.Range("E16").Value = -100
$ws.Range("F17").Copy($ws.Range("G16"))
$ws.Range("G16").Value = 2
$ws.Range("H17").Copy($ws.Range("H16"))
$ws.Range("H16").Value = -100
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = 62.5
$ws.Range("N16").Value = -87

# --- Row 17 (Fel. Assault): only M17 changes ---
$ws.Range("M17").Value = 150

# --- Row 18 (Burglary): C18/D18/E18 go from numbers back to "N/A"/"***.*" text ---
$ws.Range("C15").Copy($ws.Range("C18"))
$ws.Range("C15").Copy($ws.Range("D18"))
$ws.Range("E15").Copy($ws.Range("E18"))
$ws.Range("N18").Value = -88.235294117647

# --- Row 19 (Gr. Larceny): C19/D19/E19/F19 swap which side is "N/A" ---
$ws.Range("F17").Copy($ws.Range("C19"))
$ws.Range("C19").Value = 1
$ws.Range("C15").Copy($ws.Range("D19"))
$ws.Range("E15").Copy($ws.Range("E19"))
$ws.Range("F17").Copy($ws.Range("F19"))
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = -50
$ws.Range("I19").Value = 26
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = -21.212121212121
$ws.Range("M19").Value = -58.064516129032
$ws.Range("N19").Value = -85.227272727272

# --- Row 21 (TOTAL): value-only updates ---
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 6
$ws.Range("H21").Value = -33.333333333333
$ws.Range("I21").Value = 75
$ws.Range("J21").Value = 60
$ws.Range("K21").Value = 25
$ws.Range("L21").Value = 31.578947368421
$ws.Range("M21").Value = -25
$ws.Range("N21").Value = -84.076433121019

# --- Row 24 (Petit Larceny): C24/D24/E24 go from "N/A"/"***.*" text to real numbers ---
$ws.Range("F17").Copy($ws.Range("C24"))
$ws.Range("C24").Value = 2
$ws.Range("F17").Copy($ws.Range("D24"))
$ws.Range("D24").Value = 1
$ws.Range("H17").Copy($ws.Range("E24"))
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 3
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = -25
$ws.Range("I24").Value = 32
$ws.Range("J24").Value = 34
$ws.Range("K24").Value = -5.882352941176
$ws.Range("L24").Value = -3.030303030303
$ws.Range("M24").Value = -68.316831683168

# --- Row 25 (Misd. Assault): C25 goes from number to "N/A" text ---
$ws.Range("C15").Copy($ws.Range("C25"))
$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 2
$ws.Range("M25").Value = 68.421052631578

# --- Row 27 (Other Sex Crimes): value-only updates ---
$ws.Range("G27").Value = 2
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 71.428571428571

# --- Row 30 (Hate Crimes): D30/E30 go from numbers to "N/A"/"***.*" text ---
$ws.Range("C15").Copy($ws.Range("D30"))
$ws.Range("E15").Copy($ws.Range("E30"))
